# "Update 13-Feb-2021, midday update."
# Adds several new petty-cash ledger entries (rows 28-33) on Sheet1 ("Buku KAS
# HARIAN"-style cash book) and tops up three existing entries (rows 24-26)
# with extra amounts. The running-balance column E auto-recalculates via its
# existing shared formula, and the view's active cell moves along with the
# new data-entry point.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# --- Top up three existing rows with additional amounts ---------------------
# Row 24: Wages (D) gains another 260000
$ws.Range("D24").Formula = "=60000+260000"

# Row 25: A/R (C) gains three more receipts
$ws.Range("C25").Formula = "=2000000+762500+5840000+15740000+3850000+21101000"

# Row 26: TRANSFER BCA (D) gains two more payments
$ws.Range("D26").Formula = "=762500+15740000+2072000"

# --- New rows 28-33: continuation of the ledger ------------------------------
# Row 28: SALES - cash/retail
$ws.Range("B28").Value = "SALES - cash/retail"
$ws.Range("C28").Formula = "=40460525-9016525-21101000"

# Row 29: SELISIH - lebih
$ws.Range("B29").Value = "SELISIH - lebih"
$ws.Range("C29").Value = 10000

# Row 30: SETOR KE BANK
$ws.Range("B30").Value = "SETOR KE BANK"
$ws.Range("D30").Value = 38000000

# Row 31: new day (13-Feb-2021 serial date), Wages Expense
$ws.Range("A31").Value = 44513
$ws.Range("B31").Value = "Wages Expense"
$ws.Range("D31").Formula = "=60000"

# Row 32: A/R
$ws.Range("B32").Value = "A/R"
$ws.Range("C32").Formula = "=5250000+90000000"

# Row 33: TRANSFER BCA
$ws.Range("B33").Value = "TRANSFER BCA"
$ws.Range("D33").Formula = "=5250000+50000000+2500000"

# --- Move the view's active cell to the new data-entry point ----------------
$ws.Range("D51").Select()
